# Auto-generated Excel COM-interop script applying the Malboro_Profits market-data refresh.
# Updates columns H-N (price/profit figures) for the specific rows touched by the commit,
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 1338293.9
$ws.Range("J17").Value = 1338293.9
$ws.Range("L17").Value = 4014881.7
$ws.Range("N17").Value = -4015217.7

# Row 32
$ws.Range("H32").Value = 6586.5835
$ws.Range("I32").Value = 2528
$ws.Range("J32").Value = 9485.571
$ws.Range("K32").Value = 2528
$ws.Range("L32").Value = 9485.571
$ws.Range("M32").Value = -2202
$ws.Range("N32").Value = -10137.571

# Row 41
$ws.Range("H41").Value = 171.86667
$ws.Range("I41").Value = 117.7
$ws.Range("K41").Value = 117.7
$ws.Range("M41").Value = 322.3

# Row 45
$ws.Range("H45").Value = 2625
$ws.Range("J45").Value = 2625
$ws.Range("L45").Value = 7875
$ws.Range("N45").Value = -8259

# Row 101
$ws.Range("H101").Value = 1341.4546
$ws.Range("I101").Value = 1507
$ws.Range("K101").Value = 4521
$ws.Range("M101").Value = -2899

# Row 112
$ws.Range("H112").Value = 3487.9092
$ws.Range("J112").Value = 3931.889
$ws.Range("L112").Value = 11795.667
$ws.Range("N112").Value = -14011.667

# Row 132
$ws.Range("H132").Value = 8448.378000000001
$ws.Range("I132").Value = 6991.8945
$ws.Range("K132").Value = 20975.6835
$ws.Range("M132").Value = -18445.6835

# Row 138
$ws.Range("H138").Value = 3218.141
$ws.Range("I138").Value = 3612.3794
$ws.Range("J138").Value = 2984.8164
$ws.Range("K138").Value = 10837.1382
$ws.Range("L138").Value = 8954.449200000001
$ws.Range("M138").Value = -5697.138199999999
$ws.Range("N138").Value = -19234.4492

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 1341.2307
$ws.Range("I2").Value = 952
$ws.Range("J2").Value = 1584.5
$ws.Range("K2").Value = 952
$ws.Range("L2").Value = 1584.5
$ws.Range("M2").Value = -839
$ws.Range("N2").Value = -1810.5

# Row 4
$ws.Range("H4").Value = 817.8
$ws.Range("I4").Value = 796.3333
$ws.Range("J4").Value = 850
$ws.Range("K4").Value = 796.3333
$ws.Range("L4").Value = 850
$ws.Range("M4").Value = -680.3333
$ws.Range("N4").Value = -1082

# Row 32
$ws.Range("H32").Value = 27540.96
$ws.Range("I32").Value = 4929.7
$ws.Range("J32").Value = 61457.85
$ws.Range("K32").Value = 4929.7
$ws.Range("L32").Value = 61457.85
$ws.Range("M32").Value = -4642.7
$ws.Range("N32").Value = -62031.85

# Row 74
$ws.Range("H74").Value = 14916.879
$ws.Range("I74").Value = 2716.5557
$ws.Range("J74").Value = 29557.268
$ws.Range("K74").Value = 2716.5557
$ws.Range("L74").Value = 29557.268
$ws.Range("M74").Value = -1842.5557
$ws.Range("N74").Value = -31305.268

# Row 77
$ws.Range("H77").Value = 14916.879
$ws.Range("I77").Value = 2716.5557
$ws.Range("J77").Value = 29557.268
$ws.Range("K77").Value = 13582.7785
$ws.Range("L77").Value = 147786.34
$ws.Range("M77").Value = -9214.7785
$ws.Range("N77").Value = -156522.34

# Row 116
$ws.Range("H116").Value = 1341.2307
$ws.Range("I116").Value = 952
$ws.Range("J116").Value = 1584.5
$ws.Range("K116").Value = 952
$ws.Range("L116").Value = 1584.5
$ws.Range("M116").Value = 1342
$ws.Range("N116").Value = -6172.5

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 1341.2307
$ws.Range("I3").Value = 952
$ws.Range("J3").Value = 1584.5
$ws.Range("K3").Value = 952
$ws.Range("L3").Value = 1584.5
$ws.Range("M3").Value = -838
$ws.Range("N3").Value = -1812.5

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()

# Row 86
$ws.Range("H86").Value = 5935.3125
$ws.Range("I86").Value = 5376.222
$ws.Range("J86").Value = 6654.143
$ws.Range("K86").Value = 5376.222
$ws.Range("L86").Value = 6654.143
$ws.Range("M86").Value = -4253.222
$ws.Range("N86").Value = -8900.143

# Row 89
$ws.Range("H89").Value = 5935.3125
$ws.Range("I89").Value = 5376.222
$ws.Range("J89").Value = 6654.143
$ws.Range("K89").Value = 26881.11
$ws.Range("L89").Value = 33270.715
$ws.Range("M89").Value = -21265.11
$ws.Range("N89").Value = -44502.715

# Row 99
$ws.Range("H99").Value = 2868.2856
$ws.Range("I99").Value = 3274.5
$ws.Range("K99").Value = 3274.5
$ws.Range("M99").Value = -1776.5

# Row 107
$ws.Range("H107").Value = 2443.6216
$ws.Range("I107").Value = 3128.96
$ws.Range("J107").Value = 1015.8333
$ws.Range("K107").Value = 3128.96
$ws.Range("L107").Value = 1015.8333
$ws.Range("M107").Value = -1208.96
$ws.Range("N107").Value = -4855.8333

# Row 134
$ws.Range("H134").Value = 8473.927
$ws.Range("I134").Value = 2419.5652
$ws.Range("K134").Value = 7258.6956
$ws.Range("M134").Value = -4723.6956

# Row 138
$ws.Range("H138").Value = 78309.55
$ws.Range("J138").Value = 82140.5
$ws.Range("L138").Value = 82140.5
$ws.Range("N138").Value = -92420.5

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 22354.709
$ws.Range("I31").Value = 3675.3333
$ws.Range("J31").Value = 33562.332
$ws.Range("K31").Value = 3675.3333
$ws.Range("L31").Value = 33562.332
$ws.Range("M31").Value = -3380.3333
$ws.Range("N31").Value = -34152.332

# Row 34
$ws.Range("H34").Value = 22354.709
$ws.Range("I34").Value = 3675.3333
$ws.Range("J34").Value = 33562.332
$ws.Range("K34").Value = 3675.3333
$ws.Range("L34").Value = 33562.332
$ws.Range("M34").Value = -3473.3333
$ws.Range("N34").Value = -33966.332

# Row 100
$ws.Range("H100").Value = 60250
$ws.Range("J100").Value = 60250
$ws.Range("L100").Value = 60250
$ws.Range("N100").Value = -62414

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")

# Row 51
$ws.Range("H51").Value = 171768720
$ws.Range("I51").Value = 10204081
$ws.Range("J51").Value = 333333340
$ws.Range("K51").Value = 30612243
$ws.Range("L51").Value = 1000000020
$ws.Range("M51").Value = -30611783
$ws.Range("N51").Value = -1000000940

# Row 75
$ws.Range("H75").Value = 1742.3334
$ws.Range("J75").Value = 2851.4
$ws.Range("L75").Value = 8554.200000000001
$ws.Range("N75").Value = -10550.2

# Row 78
$ws.Range("H78").Value = 1742.3334
$ws.Range("J78").Value = 2851.4
$ws.Range("L78").Value = 25662.6
$ws.Range("N78").Value = -35646.60000000001

# Row 101
$ws.Range("H101").Value = 10287.556
$ws.Range("J101").Value = 10287.556
$ws.Range("L101").Value = 30862.668
$ws.Range("N101").Value = -35730.66800000001

# Row 112
$ws.Range("H112").Value = 13466.667
$ws.Range("I112").Value = 5000
$ws.Range("J112").Value = 14071.429
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 42214.287
$ws.Range("M112").Value = -13892
$ws.Range("N112").Value = -44430.287

# Row 122
$ws.Range("H122").Value = 9357931
$ws.Range("I122").Value = 23359124
$ws.Range("K122").Value = 210232116
$ws.Range("M122").Value = -210229666

# Row 131
$ws.Range("H131").Value = 1476.2
$ws.Range("J131").Value = 1481.1531
$ws.Range("L131").Value = 4443.4593
$ws.Range("N131").Value = -14523.4593

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")

# Row 52
$ws.Range("H52").Value = 47977.5
$ws.Range("J52").Value = 47977.5
$ws.Range("L52").Value = 47977.5
$ws.Range("N52").Value = -48495.5

# Row 80
$ws.Range("H80").Value = 19076.23
$ws.Range("I80").Value = 12124.5
$ws.Range("K80").Value = 12124.5
$ws.Range("M80").Value = -11126.5

# Row 83
$ws.Range("H83").Value = 19076.23
$ws.Range("I83").Value = 12124.5
$ws.Range("K83").Value = 60622.5
$ws.Range("M83").Value = -55630.5

# Row 132
$ws.Range("H132").Value = 6498.5415
$ws.Range("I132").Value = 2576.4736
$ws.Range("K132").Value = 7729.4208
$ws.Range("M132").Value = -5199.4208

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")

# Row 61
$ws.Range("H61").Value = 3453.48
$ws.Range("I61").Value = 1341.9286
$ws.Range("K61").Value = 1341.9286
$ws.Range("M61").Value = -1139.9286

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 113
$ws.Range("H113").Value = 3453.48
$ws.Range("I113").Value = 1341.9286
$ws.Range("K113").Value = 1341.9286
$ws.Range("M113").Value = 828.0714

# Row 132
$ws.Range("H132").Value = 1913282.1
$ws.Range("I132").Value = 824.7778
$ws.Range("K132").Value = 2474.3334
$ws.Range("M132").Value = 55.66660000000002

# Row 136
$ws.Range("H136").Value = 11663.653
$ws.Range("I136").Value = 12311.423
$ws.Range("J136").Value = 10931.392
$ws.Range("K136").Value = 36934.269
$ws.Range("L136").Value = 32794.176
$ws.Range("M136").Value = -34384.269
$ws.Range("N136").Value = -37894.176

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")

# Row 69
$ws.Range("H69").Value = 45817.375
$ws.Range("J69").Value = 45817.375
$ws.Range("L69").Value = 45817.375
$ws.Range("N69").Value = -47315.375

# Row 72
$ws.Range("H72").Value = 45817.375
$ws.Range("J72").Value = 45817.375
$ws.Range("L72").Value = 137452.125
$ws.Range("N72").Value = -144940.125

# Row 127
$ws.Range("H127").Value = 24666.666
$ws.Range("J127").Value = 24666.666
$ws.Range("L127").Value = 24666.666
$ws.Range("N127").Value = -34586.666

# Row 136
$ws.Range("H136").Value = 9607.870999999999
$ws.Range("I136").Value = 1365.2858
$ws.Range("K136").Value = 4095.8574
$ws.Range("M136").Value = -1545.8574

